$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1 data edit: B18 18 -> 20 (dependent formulas E18:K18 / M18:S18 recalc automatically) ---
$ws1.Range("B18").Value = 20

# --- Sheet1 row 19: newly entered forecast-threshold annotations ---
$ws1.Range("E19").Value = "大于 10.0"
$ws1.Range("F19").Value = 40
$ws1.Range("G19").Value = "大于 63.5"
$ws1.Range("H19").Value = "大于 11.0"

# Apply the same number format as the rest of the table (style s="4") to row 19's
# remaining cells (I19:K19, M19:S19) and to the already-valued cells, matching the
# blank-but-styled cells added alongside the new values.
$ws1.Range("E19:K19").NumberFormat = "0.0"
$ws1.Range("M19:S19").NumberFormat = "0.0"

# --- Sheet1 rows 20:82: same style fill (blank styled cells, no values) ---
$ws1.Range("E20:K82").NumberFormat = "0.0"
$ws1.Range("M20:S82").NumberFormat = "0.0"

# --- Sheet1 styles: right-align the "0.0" formatted forecast/error columns ---
$ws1.Range("E2:K82").HorizontalAlignment = -4152
$ws1.Range("M2:S82").HorizontalAlignment = -4152

# --- Sheet1 view: selection / active cell moves to I19, Sheet1 becomes the active tab
#     (Chart1 loses tabSelected / workbook activeTab as a side effect) ---
$ws1.Activate()
$ws1.Range("I19").Select()
